$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- Edit 1: Title shape - append "(2022? ~)" to "Global Season Stage" ---
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleRun = $titleTr.Characters(14, 21)
$titleRun.Text = "– Global Season Stage (2022? ~)"

# --- Edit 2: Content shape - rework the "6천만 달러..." sentence ---
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange

# Replace the trailing three runs (old text: "천만 달러 이상의 돈을 내야지 운영에
# 참여하는 것도 중요하고" + ", " + "어느 정도의 액수를 내야 운영이 참여하게끔 그런
# 제도도 하나씩 만드는 것이다") with a single merged run of new text.
# Do this BEFORE editing the leading "6" run so offsets stay valid.
$mergedRun = $bodyTr.Characters(513, 77)
$mergedRun.Text = "천만 달러 혹은 그 이상의 돈을 내야지 운영에 참여하는 것등 여러가지 리그 운영에 대한 장치를 마련해야 한다"

# Replace "6" with "3~6"
$numRun = $bodyTr.Characters(512, 1)
$numRun.Text = "3~6"
